$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.375.19"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.509.01"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +5.93%  "
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("D12").Value = "4.105.19"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "3.505.60"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "64.367.26"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.579"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "3.648.46"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.85%  "
$ws.Range("E33").Value = "  +5.99%  "
$ws.Range("D34").Value = "3.538.86"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").Value = "2.376.09"
$ws.Range("E49").Value = "  -4.48%  "
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("E51").Value = "  -0.12%  "
